$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 121.8
$ws.Range("I9").Value = 136.66667
$ws.Range("J9").Value = 99.5
$ws.Range("K9").Value = 136.66667
$ws.Range("L9").Value = 99.5
$ws.Range("M9").Value = 32.33332999999999
$ws.Range("N9").Value = -437.5
$ws.Range("H29").Value = 1837.5
$ws.Range("H58").Value = 510.77777
$ws.Range("J58").Value = 2017
$ws.Range("L58").Value = 6051
$ws.Range("N58").Value = -6351
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").ClearContents() | Out-Null
$ws.Range("M92").ClearContents() | Out-Null
$ws.Range("N92").Value = 0
$ws.Range("H129").Value = 855.4032
$ws.Range("I129").Value = 472.85715
$ws.Range("J129").Value = 966.9792
$ws.Range("K129").Value = 1418.57145
$ws.Range("L129").Value = 2900.9376
$ws.Range("M129").Value = 3581.42855
$ws.Range("N129").Value = -12900.9376
$ws.Range("H137").Value = 3574269
$ws.Range("I137").Value = 8336025
$ws.Range("J137").Value = 2951.75
$ws.Range("K137").Value = 25008075
$ws.Range("L137").Value = 8855.25
$ws.Range("M137").Value = -25005525
$ws.Range("N137").Value = -13955.25

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10789897
$ws.Range("I32").Value = 13934854
$ws.Range("K32").Value = 13934854
$ws.Range("M32").Value = -13934567
$ws.Range("H61").Value = 45547130
$ws.Range("I61").Value = 52685972
$ws.Range("J61").Value = 334433.34
$ws.Range("K61").Value = 52685972
$ws.Range("L61").Value = 334433.34
$ws.Range("M61").Value = -52685760
$ws.Range("N61").Value = -334857.34
$ws.Range("H74").Value = 18001076
$ws.Range("I74").Value = 27889692
$ws.Range("J74").Value = 201570
$ws.Range("K74").Value = 27889692
$ws.Range("L74").Value = 201570
$ws.Range("M74").Value = -27888818
$ws.Range("N74").Value = -203318
$ws.Range("H77").Value = 18001076
$ws.Range("I77").Value = 27889692
$ws.Range("J77").Value = 201570
$ws.Range("K77").Value = 139448460
$ws.Range("L77").Value = 1007850
$ws.Range("M77").Value = -139444092
$ws.Range("N77").Value = -1016586
$ws.Range("H119").Value = 42698
$ws.Range("J119").Value = 42698
$ws.Range("L119").Value = 42698
$ws.Range("N119").Value = -52374
$ws.Range("H122").Value = 1303.36
$ws.Range("I122").Value = 1295.1666
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3885.4998
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1435.4998
$ws.Range("N122").Value = -9400
$ws.Range("H136").Value = 45547130
$ws.Range("I136").Value = 52685972
$ws.Range("J136").Value = 334433.34
$ws.Range("K136").Value = 158057916
$ws.Range("L136").Value = 1003300.02
$ws.Range("M136").Value = -158055366
$ws.Range("N136").Value = -1008400.02

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 568.5714
$ws.Range("I37").Value = 330
$ws.Range("K37").Value = 330
$ws.Range("M37").Value = -193
$ws.Range("H64").Value = 965
$ws.Range("H67").Value = 965

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1664.6522
$ws.Range("I31").Value = 1513.3182
$ws.Range("J31").Value = 4994
$ws.Range("K31").Value = 1513.3182
$ws.Range("L31").Value = 4994
$ws.Range("M31").Value = -1218.3182
$ws.Range("N31").Value = -5584
$ws.Range("H34").Value = 1664.6522
$ws.Range("I34").Value = 1513.3182
$ws.Range("J34").Value = 4994
$ws.Range("K34").Value = 1513.3182
$ws.Range("L34").Value = 4994
$ws.Range("M34").Value = -1311.3182
$ws.Range("N34").Value = -5398
$ws.Range("H123").Value = 44780
$ws.Range("J123").Value = 44780
$ws.Range("L123").Value = 44780
$ws.Range("N123").Value = -54580
$ws.Range("H132").Value = 130343.5
$ws.Range("I132").Value = 6187.5
$ws.Range("J132").Value = 254499.5
$ws.Range("K132").Value = 18562.5
$ws.Range("L132").Value = 763498.5
$ws.Range("M132").Value = -16032.5
$ws.Range("N132").Value = -768558.5
$ws.Range("H134").Value = 46256.832
$ws.Range("I134").Value = 1210.25
$ws.Range("K134").Value = 3630.75
$ws.Range("M134").Value = -1095.75

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 42530.793
$ws.Range("J5").Value = 1011.1539
$ws.Range("L5").Value = 3033.4617
$ws.Range("N5").Value = -3257.4617
$ws.Range("H39").Value = 1920
$ws.Range("J39").Value = 1920
$ws.Range("L39").Value = 5760
$ws.Range("N39").Value = -6348
$ws.Range("H55").Value = 3112.25
$ws.Range("I55").Value = 1800
$ws.Range("J55").Value = 3199.7334
$ws.Range("K55").Value = 5400
$ws.Range("L55").Value = 9599.200199999999
$ws.Range("M55").Value = -5223
$ws.Range("N55").Value = -9953.200199999999
$ws.Range("H117").Value = 9525113
$ws.Range("I117").Value = 1043
$ws.Range("J117").Value = 16668166
$ws.Range("K117").Value = 3129
$ws.Range("L117").Value = 50004498
$ws.Range("M117").Value = 313
$ws.Range("N117").Value = -50011382
$ws.Range("H122").Value = 1212.625
$ws.Range("I122").Value = 516.5
$ws.Range("J122").Value = 1444.6666
$ws.Range("K122").Value = 4648.5
$ws.Range("L122").Value = 13001.9994
$ws.Range("M122").Value = -2198.5
$ws.Range("N122").Value = -17901.9994
$ws.Range("H132").Value = 2190.6
$ws.Range("I132").Value = 4243.5
$ws.Range("K132").Value = 38191.5
$ws.Range("M132").Value = -35661.5
$ws.Range("H135").Value = 42530.793
$ws.Range("J135").Value = 1011.1539
$ws.Range("L135").Value = 9100.3851
$ws.Range("N135").Value = -14170.3851

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents() | Out-Null
$ws.Range("H93").Value = 22184.066
$ws.Range("H99").Value = 5510
$ws.Range("I99").Value = 5510
$ws.Range("K99").Value = 5510
$ws.Range("M99").Value = -3264
$ws.Range("H122").Value = 3428.4285
$ws.Range("I122").Value = 1333.3334
$ws.Range("J122").Value = 4999.75
$ws.Range("K122").Value = 4000.0002
$ws.Range("L122").Value = 14999.25
$ws.Range("M122").Value = -1550.0002
$ws.Range("N122").Value = -19899.25
$ws.Range("H123").Value = 21376.818
$ws.Range("J123").Value = 21376.818
$ws.Range("L123").Value = 21376.818
$ws.Range("N123").Value = -26276.818
$ws.Range("H132").Value = 55841.5
$ws.Range("I132").Value = 41299.5
$ws.Range("J132").Value = 87349.164
$ws.Range("K132").Value = 123898.5
$ws.Range("L132").Value = 262047.492
$ws.Range("M132").Value = -121368.5
$ws.Range("N132").Value = -267107.492

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 54227.895
$ws.Range("I100").Value = 201380
$ws.Range("J100").Value = 1673.5714
$ws.Range("K100").Value = 201380
$ws.Range("L100").Value = 1673.5714
$ws.Range("M100").Value = -200839
$ws.Range("N100").Value = -2755.5714
$ws.Range("H122").Value = 3453.182
$ws.Range("I122").Value = 3708.5715
$ws.Range("J122").Value = 3006.25
$ws.Range("K122").Value = 11125.7145
$ws.Range("L122").Value = 9018.75
$ws.Range("M122").Value = -8675.7145
$ws.Range("N122").Value = -13918.75
$ws.Range("H132").Value = 66387.56
$ws.Range("I132").Value = 3355.3333
$ws.Range("J132").Value = 147429
$ws.Range("K132").Value = 10065.9999
$ws.Range("L132").Value = 442287
$ws.Range("M132").Value = -7535.999899999999
$ws.Range("N132").Value = -447347
$ws.Range("H136").Value = 127643.375
$ws.Range("I136").Value = 74256.71000000001
$ws.Range("J136").Value = 501350
$ws.Range("K136").Value = 222770.13
$ws.Range("L136").Value = 1504050
$ws.Range("M136").Value = -220220.13
$ws.Range("N136").Value = -1509150

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 44600
$ws.Range("J46").Value = 44600
$ws.Range("L46").Value = 44600
$ws.Range("N46").Value = -45062
$ws.Range("H54").Value = 13638.556
$ws.Range("J54").Value = 13638.556
$ws.Range("L54").Value = 13638.556
$ws.Range("N54").Value = -14678.556
$ws.Range("H62").Value = 4002.3
$ws.Range("I62").Value = 4000.6667
$ws.Range("K62").Value = 4000.6667
$ws.Range("M62").Value = -3376.6667
$ws.Range("H65").Value = 4002.3
$ws.Range("I65").Value = 4000.6667
$ws.Range("K65").Value = 20003.3335
$ws.Range("M65").Value = -16883.3335
$ws.Range("H81").Value = 3987.4285
$ws.Range("I81").Value = 3120.4
$ws.Range("J81").Value = 4469.1113
$ws.Range("K81").Value = 6240.8
$ws.Range("L81").Value = 8938.222599999999
$ws.Range("M81").Value = -5179.8
$ws.Range("N81").Value = -11060.2226
$ws.Range("H84").Value = 3987.4285
$ws.Range("I84").Value = 3120.4
$ws.Range("J84").Value = 4469.1113
$ws.Range("K84").Value = 31204
$ws.Range("L84").Value = 44691.113
$ws.Range("M84").Value = -25900
$ws.Range("N84").Value = -55299.113
$ws.Range("H115").Value = 29929.1
$ws.Range("J115").Value = 29929.1
$ws.Range("L115").Value = 29929.1
$ws.Range("N115").Value = -33063.1
$ws.Range("H122").Value = 2298.9
$ws.Range("I122").Value = 2048.4285
$ws.Range("J122").Value = 2883.3333
$ws.Range("K122").Value = 6145.2855
$ws.Range("L122").Value = 8649.999899999999
$ws.Range("M122").Value = -3695.2855
$ws.Range("N122").Value = -13549.9999
$ws.Range("H132").Value = 252124.75
$ws.Range("I132").Value = 334666.66
$ws.Range("J132").Value = 202599.6
$ws.Range("K132").Value = 1003999.98
$ws.Range("L132").Value = 607798.8
$ws.Range("M132").Value = -1001469.98
$ws.Range("N132").Value = -612858.8
$ws.Range("H134").Value = 44600
$ws.Range("J134").Value = 44600
$ws.Range("L134").Value = 133800
$ws.Range("N134").Value = -138870

